# Swap the order of "System" and the email address in column G ("Recorded By")
# wherever the value is exactly "System, dnasr281@gmail.com", turning it into
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$startRow = $used.Row
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = $startRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ($text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
